$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with the latest values.
# The D column is forced to Text format before assignment so that values such as
# "1.00" or "0.0787" are preserved verbatim as strings instead of being parsed as numbers.

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = '42.853.64'
$dCell.Style = "Normal"
$ws.Range("E2").Value = '  -0.58%  '

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = '2.296.38'
$dCell.Style = "Normal"
$ws.Range("E3").Value = '  -1.00%  '

$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = '1.00'
$dCell.Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '299.60'
$dCell.Style = "Normal"
$ws.Range("E5").Value = '  -1.26%  '

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = '97.32'
$dCell.Style = "Normal"
$ws.Range("E6").Value = '  -2.42%  '

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = '0.514'
$dCell.Style = "Normal"
$ws.Range("E7").Value = '  +1.16%  '

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = '1.00'
$dCell.Style = "Normal"
$ws.Range("E8").Value = '  -0.08%  '

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = '0.506'
$dCell.Style = "Normal"
$ws.Range("E9").Value = '  -2.51%  '

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = '36.03'
$dCell.Style = "Normal"
$ws.Range("E10").Value = '  -0.19%  '

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0787'
$dCell.Style = "Normal"
$ws.Range("E11").Value = '  -0.44%  '

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = '0.117'
$dCell.Style = "Normal"
$ws.Range("E12").Value = '  +0.75%  '

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = '17.71'
$dCell.Style = "Normal"
$ws.Range("E13").Value = '  -0.18%  '

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = '6.76'
$dCell.Style = "Normal"
$ws.Range("E14").Value = '  -2.29%  '

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = '2.644.87'
$dCell.Style = "Normal"
$ws.Range("E15").Value = '  -1.39%  '

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = '2.300.99'
$dCell.Style = "Normal"
$ws.Range("E16").Value = '  -2.03%  '

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = '0.777'
$dCell.Style = "Normal"
$ws.Range("E17").Value = '  -2.27%  '

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = '42.811.38'
$dCell.Style = "Normal"
$ws.Range("E18").Value = '  -0.54%  '

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = '12.52'
$dCell.Style = "Normal"
$ws.Range("E19").Value = '  -4.93%  '

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0₃0907'
$dCell.Style = "Normal"
$ws.Range("E20").Value = '  -0.23%  '

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = '6.06'
$dCell.Style = "Normal"
$ws.Range("E21").Value = '  -2.15%  '

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = '67.92'
$dCell.Style = "Normal"
$ws.Range("E22").Value = '  -0.53%  '

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '241.71'
$dCell.Style = "Normal"
$ws.Range("E23").Value = '  +0.73%  '

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = '2.14'
$dCell.Style = "Normal"
$ws.Range("E24").Value = '  -1.09%  '

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = '1.00'
$dCell.Style = "Normal"
$ws.Range("E25").Value = '  +0.02%  '

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = '2.42'
$dCell.Style = "Normal"
$ws.Range("E26").Value = '  -1.50%  '

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = '4.02'
$dCell.Style = "Normal"
$ws.Range("E27").Value = '  -0.33%  '

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = '25.19'
$dCell.Style = "Normal"
$ws.Range("E28").Value = '  -1.35%  '

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = '165.98'
$dCell.Style = "Normal"
$ws.Range("E29").Value = '  -1.92%  '

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = '2.03'
$dCell.Style = "Normal"
$ws.Range("E30").Value = '  -1.03%  '

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = '9.01'
$dCell.Style = "Normal"
$ws.Range("E31").Value = '  -2.00%  '

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = '32.82'
$dCell.Style = "Normal"
$ws.Range("E32").Value = '  -3.76%  '

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = '1.00'
$dCell.Style = "Normal"
$ws.Range("E33").Value = '  +0.02%  '

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = '4.77'
$dCell.Style = "Normal"
$ws.Range("E34").Value = '  -4.11%  '

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = '4.99'
$dCell.Style = "Normal"
$ws.Range("E35").Value = '  -3.50%  '

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = '17.15'
$dCell.Style = "Normal"
$ws.Range("E36").Value = '  -4.32%  '

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = '2.38'
$dCell.Style = "Normal"
$ws.Range("E37").Value = '  -0.62%  '

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0688'
$dCell.Style = "Normal"
$ws.Range("E38").Value = '  -1.46%  '

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = '0.100'
$dCell.Style = "Normal"
$ws.Range("E39").Value = '  -2.26%  '

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = '1.76'
$dCell.Style = "Normal"
$ws.Range("E40").Value = '  -3.68%  '

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = '2.74'
$dCell.Style = "Normal"
$ws.Range("E41").Value = '  -1.38%  '

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = '0.110'
$dCell.Style = "Normal"
$ws.Range("E42").Value = '  -0.08%  '

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = '2.012.03'
$dCell.Style = "Normal"
$ws.Range("E43").Value = '  +0.97%  '

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0284'
$dCell.Style = "Normal"
$ws.Range("E44").Value = '  -1.81%  '

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = '10.11'
$dCell.Style = "Normal"
$ws.Range("E45").Value = '  -0.83%  '

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = '2.12'
$dCell.Style = "Normal"
$ws.Range("E46").Value = '  -5.37%  '

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = '17.19'
$dCell.Style = "Normal"
$ws.Range("E47").Value = '  -1.85%  '

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = '2.77'
$dCell.Style = "Normal"
$ws.Range("E48").Value = '  -2.89%  '

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = '2.523.66'
$dCell.Style = "Normal"
$ws.Range("E49").Value = '  -0.94%  '

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = '53.15'
$dCell.Style = "Normal"
$ws.Range("E50").Value = '  -3.44%  '

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = '2.77'
$dCell.Style = "Normal"
$ws.Range("E51").Value = '  -8.46%  '
